$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.4145870963267839,0.1615748415687079,0.1784203654063603,1.907250001524517,1.258863221492163,1.139205889773933,0,0.2336961426485047,0,0.1609786650762786,0,2.873581687696401,4.892343549182669),
    @(0.4107024072645231,0.1615709920218293,0.1778979785888062,1.877333202784783,1.226370776659223,1.129112780806395,0,0.2324819634783069,0,0.1611842094933245,0,2.562605684679681,4.801883612414144),
    @(0.4085160041096572,0.1616226682651245,0.1776499131431954,1.859969350165471,1.207147119790932,1.123460027337075,0,0.231840510704231,0,0.1613630788569438,0,2.371325805375818,4.749003842339505),
    @(0.4076750433124943,0.1616573420264942,0.1775671208319913,1.853146138554038,1.199495870671711,1.121293439370447,0,0.2316052930081511,0,0.1614492342834453,0,2.293303068607429,4.728124154461909),
    @(0.4075384239871767,0.1616639226908205,0.1775544788788963,1.852028410230858,1.198236405631889,1.120941952292355,0,0.2315678168707933,0,0.1614643421285713,0,2.280343261403573,4.7246975083238),
    @(0.4085044600711711,0.1616230807239916,0.1776487224689518,1.859876306866127,1.207043193539619,1.123430253393025,0,0.2318372324579556,0,0.1613641870416771,0,2.370273851395496,4.748719542336403),
    @(0.4132064037997338,0.1615623138236018,0.1782251725285136,1.896726010695033,1.247508824146593,1.135612805412876,0,0.2332559032987476,0,0.1610386191533593,0,2.766433886209825,4.86059998215768),
    @(0.4240043641161151,0.1618709033210948,0.1799316716272514,1.976973078055465,1.332644793882565,1.163823441513983,0,0.2368634153445655,0,0.1608171973377601,0,3.540180268007646,5.101162523401683),
    @(0.4329007671987029,0.1623572111729956,0.1815361776646149,2.040820501381816,1.398750128412559,1.187188775875882,0,0.2400175384909247,0,0.1609077591742292,0,4.10623028343673,5.290882596801623),
    @(0.4371575361408873,0.1626345781302803,0.1823422034696165,2.070933559428283,1.429602432515935,1.198392676900625,0,0.2415619452896181,0,0.1610037532333237,0,4.363110593465422,5.38002757159876),
    @(0.4387996321026719,0.1627476604520766,0.1826583594944289,2.08249050447148,1.441398060645895,1.202718009386928,0,0.2421625299244141,0,0.1610479651501322,0,4.460285735714251,5.414193847512422),
    @(0.4384446366101997,0.1627229485253068,0.1825897836986528,2.079994667914889,1.438852651878591,1.201782795529283,0,0.2420324828273266,0,0.1610380939957068,0,4.439361943450422,5.406817328695752),
    @(0.4372920281793142,0.1626437202922801,0.1823679948194759,2.071881273450145,1.430570609092314,1.198746868062074,0,0.2416110401029954,0,0.1610072331331551,0,4.371107314139522,5.382830251232349),
    @(0.4365899483082671,0.162596238288593,0.1822335658805763,2.066931616112015,1.425512281660389,1.196898039116945,0,0.2413549453743187,0,0.1609893531537594,0,4.329286057409945,5.368190755786713),
    @(0.4326267970158142,0.162340212290303,0.1814850325341908,2.03887405489867,1.396749589797167,1.186468141939287,0,0.2399188123669802,0,0.1609025867895397,0,4.089429168003562,5.285113966609799),
    @(0.4302492454540641,0.1621975111976752,0.1810453185918668,2.021935419282357,1.379304747435498,1.18021697327589,0,0.2390658528425007,0,0.1608633798883545,0,3.94211849063862,5.234876931550275),
    @(0.4289014840593381,0.1621207191674614,0.1807995731982359,2.012293332637199,1.369344397988556,1.176675578534685,0,0.2385855676477391,0,0.1608459879753426,0,3.857331695637754,5.206249231189418),
    @(0.4284485465293244,0.1620956272874636,0.1807175993518335,2.009045956693356,1.365984603809437,1.175485817301279,0,0.2384247230672329,0,0.1608409859849402,0,3.828614786364199,5.196602277993406),
    @(0.4305002966017071,0.162212155079132,0.1810913853069387,2.023728156611384,1.381154175649726,1.180876820556875,0,0.2391555843695983,0,0.1608670197445861,0,3.957806003280837,5.240197079765721),
    @(0.4376297590455067,0.1626667732736564,0.1824328430329878,2.074260199388647,1.43300018958746,1.199636350297226,0,0.2417344006085784,0,0.1610160845176338,0,4.391158149571254,5.389864732357523),
    @(0.4424650069452412,0.1630108001389914,0.1833732736285576,2.108182241695829,1.467540598247496,1.2123785806412,0,0.2435116142159544,0,0.1611593186794593,0,4.673791817957863,5.49006548527683),
    @(0.4398682683589925,0.1628229022956802,0.1828655236690651,2.089995342099598,1.449045613039544,1.205533731124149,0,0.2425546833317895,0,0.1610786862503701,0,4.523002190005627,5.436368103370228),
    @(0.4303867368071792,0.1622055182295554,0.1810705365783676,2.022917360543119,1.380317834488125,1.180578340110031,0,0.2391149853095627,0,0.1608653581273956,0,3.950713976768498,5.237791048862846),
    @(0.4209142018309251,0.1617417508117356,0.1794083929518848,1.954407016463563,1.308991014268713,1.155728654418482,0,0.235799068926255,0,0.1608325559542507,0,3.331249627311365,5.033810303472592)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    $arr = New-Object 'object[,]' 1,13
    for ($j = 0; $j -lt 13; $j++) {
        $arr[0,$j] = $rowVals[$j]
    }
    $ws.Range("C" + $rowNum + ":O" + $rowNum).Value = $arr
}

Write-Host "Updated rows 2-25, columns C:O"
